$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a "Device Name" column as the new first column of the report table ---
# Only the header row (8) and the data row (9) of the position table actually
# grow a new column in the real template; the rows above (the report/device/
# period info block in rows 1-7) keep their original A:H layout. Using
# Range.Insert() would shift every row in the sheet, so instead we manually
# slide the existing H:A cells (values + formatting) one column to the right,
# working from the rightmost column inward so nothing is clobbered.
$ws.Range("H8:H9").Copy($ws.Range("I8:I9"))
$ws.Range("G8:G9").Copy($ws.Range("H8:H9"))
$ws.Range("F8:F9").Copy($ws.Range("G8:G9"))
$ws.Range("E8:E9").Copy($ws.Range("F8:F9"))
$ws.Range("D8:D9").Copy($ws.Range("E8:E9"))
$ws.Range("C8:C9").Copy($ws.Range("D8:D9"))
$ws.Range("B8:B9").Copy($ws.Range("C8:C9"))
$ws.Range("A8:A9").Copy($ws.Range("B8:B9"))

# New column A keeps the same (already-copied) header/data style, just with
# new text: the header label and the templated device-name placeholder.
$ws.Range("A8").Value = "Device Name"
$ws.Range("A9").Value = "`${position.deviceName}"

# --- Update the jxls comment hints: the table's last column moved from H to I ---
$ws.Range("A1").Comment.Text("jx:area(lastCell=""I9"")")
$ws.Range("A2").Comment.Text("jx:each(items=""devices"", var=""device"", lastCell=""I9"" multisheet=""sheetNames"")")
$ws.Range("A9").Comment.Text("jx:each(items=""device.objects"", var=""position"", lastCell=""I9"")")
